$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update D2, E2
$ws.Range("D2").Value = '66.960.43'
$ws.Range("E2").Value = '  -2.01%  '

# Row 3: update D3, E3
$ws.Range("D3").Value = '3.484.06'
$ws.Range("E3").Value = '  -2.33%  '

# Row 4: update E4
$ws.Range("E4").Value = '  +0.08%  '

# Row 5: update D5, E5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.81'
$ws.Range("E5").Value = '  -2.88%  '

# Row 6: update D6, E6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.13'
$ws.Range("E6").Value = '  -4.61%  '

# Row 7: update D7, E7
$ws.Range("D7").Value = '3.482.27'
$ws.Range("E7").Value = '  -2.32%  '

# Row 8: update E8
$ws.Range("E8").Value = '  -0.06%  '

# Row 9: update D9, E9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.480'
$ws.Range("E9").Value = '  -2.37%  '

# Row 10: update E10
$ws.Range("E10").Value = '  -2.58%  '

# Row 11: update D11, E11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.67'
$ws.Range("E11").Value = '  +3.69%  '

# Row 12: update E12
$ws.Range("E12").Value = '  -3.40%  '

# Row 13: update E13
$ws.Range("E13").Value = '  -3.79%  '

# Row 14: update D14, E14
$ws.Range("D14").Value = '4.072.49'
$ws.Range("E14").Value = '  -2.34%  '

# Row 15: update D15, E15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '31.17'
$ws.Range("E15").Value = '  -6.21%  '

# Row 16: update D16, E16
$ws.Range("D16").Value = '3.473.24'
$ws.Range("E16").Value = '  -2.49%  '

# Row 17: update D17, E17
$ws.Range("D17").Value = '66.894.10'
$ws.Range("E17").Value = '  -2.16%  '

# Row 18: update E18
$ws.Range("E18").Value = '  +0.37%  '

# Row 19: update E19
$ws.Range("E19").Value = '  -4.97%  '

# Row 20: update D20, E20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.19'
$ws.Range("E20").Value = '  +1.46%  '

# Row 21: update D21, E21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.23'
$ws.Range("E21").Value = '  -4.84%  '

# Row 22: update D22, E22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '434.72'
$ws.Range("E22").Value = '  -4.35%  '

# Row 23: update E23
$ws.Range("E23").Value = '  -5.80%  '

# Row 24: update D24, E24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '79.13'
$ws.Range("E24").Value = '  +0.82%  '

# Row 26: update D26, E26
$ws.Range("D26").Value = '3.619.87'
$ws.Range("E26").Value = '  -2.40%  '

# Row 27: update E27
$ws.Range("E27").Value = '  -10.22%  '

# Row 28: update D28, E28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.81'
$ws.Range("E28").Value = '  -7.04%  '

# Row 29: update D29, E29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.35'
$ws.Range("E29").Value = '  -9.45%  '

# Row 30: update E30
$ws.Range("E30").Value = '  -3.32%  '

# Row 31: update E31
$ws.Range("E31").Value = '  -7.60%  '

# Row 32: update B32, C32, D32, E32
$ws.Range("B32").Value = 'Binance-PegBSC-USD'
$ws.Range("C32").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.01'
$ws.Range("E32").Value = '  +0.39%  '

# Row 33: update B33, C33, D33, E33
$ws.Range("B33").Value = 'Kaspa'
$ws.Range("C33").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.167'
$ws.Range("E33").Value = '  -2.36%  '

# Row 34: update E34
$ws.Range("E34").Value = '  -3.16%  '

# Row 35: update D35, E35
$ws.Range("D35").Value = '3.473.26'
$ws.Range("E35").Value = '  -2.45%  '

# Row 36: update E36
$ws.Range("E36").Value = '  -6.24%  '

# Row 37: update D37, E37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.92'
$ws.Range("E37").Value = '  -7.29%  '

# Row 38: update E38
$ws.Range("E38").Value = '  +0.01%  '

# Row 39: update E39
$ws.Range("E39").Value = '  -4.28%  '

# Row 40: update E40
$ws.Range("E40").Value = '  +0.03%  '

# Row 41: update D41, E41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '173.87'
$ws.Range("E41").Value = '  -3.81%  '

# Row 42: update E42
$ws.Range("E42").Value = '  -3.69%  '

# Row 43: update E43
$ws.Range("E43").Value = '  -13.17%  '

# Row 44: update D44, E44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.41'
$ws.Range("E44").Value = '  -3.62%  '

# Row 45: update D45, E45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.899'
$ws.Range("E45").Value = '  +0.02%  '

# Row 46: update D46, E46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '46.46'
$ws.Range("E46").Value = '  +0.48%  '

# Row 47: update D47, E47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '28.90'
$ws.Range("E47").Value = '  -7.31%  '

# Row 48: update E48
$ws.Range("E48").Value = '  -6.98%  '

# Row 49: update D49, E49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.45'
$ws.Range("E49").Value = '  -4.36%  '

# Row 50: update D50, E50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.41'
$ws.Range("E50").Value = '  -9.65%  '

# Row 51: update D51, E51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.970'
$ws.Range("E51").Value = '  -4.90%  '
